$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data feed added a new observation (rows 771-772) for the
# "Femacal de La Calera - Brócoli" series. That pushes every existing data
# row down by two positions; the pair that previously fell off the bottom
# of the range (old rows 821-822) simply lands on new rows 823-824.
#
# Inserting two physical rows above the old row 771 reproduces exactly
# that: all formatting/styles/values of rows 771-822 slide down to
# 773-824 in one shot (Excel grows the used range to R824 automatically),
# so afterwards we only need to populate the two brand-new rows with the
# new observation's data.

$ws.Rows("771:772").Insert()

# Row 771 - "Primera" quality
$ws.Range("A771").Value2 = 3
$ws.Range("B771").Value2 = "Femacal de La Calera"
$ws.Range("C771").Value2 = "Coquimbo"
$ws.Range("D771").Value2 = 44826
$ws.Range("E771").Value2 = 5
$ws.Range("F771").Value2 = 100112023
$ws.Range("G771").Value2 = "Br" + [char]243 + "coli"
$ws.Range("H771").Value2 = "Sin especificar"
$ws.Range("I771").Value2 = "Primera"
$ws.Range("J771").Value2 = 3600
$ws.Range("K771").Value2 = 800
$ws.Range("L771").Value2 = 850
$ws.Range("M771").Value2 = 821
$ws.Range("N771").Value2 = "`$/unidad"
$ws.Range("O771").Value2 = "Provincia de Quillota"
$ws.Range("P771").Value2 = 821
$ws.Range("Q771").Value2 = 1
$ws.Range("R771").Value2 = "Hortaliza"

# Row 772 - "Segunda" quality
$ws.Range("A772").Value2 = 3
$ws.Range("B772").Value2 = "Femacal de La Calera"
$ws.Range("C772").Value2 = "Coquimbo"
$ws.Range("D772").Value2 = 44826
$ws.Range("E772").Value2 = 5
$ws.Range("F772").Value2 = 100112023
$ws.Range("G772").Value2 = "Br" + [char]243 + "coli"
$ws.Range("H772").Value2 = "Sin especificar"
$ws.Range("I772").Value2 = "Segunda"
$ws.Range("J772").Value2 = 1300
$ws.Range("K772").Value2 = 600
$ws.Range("L772").Value2 = 600
$ws.Range("M772").Value2 = 600
$ws.Range("N772").Value2 = "`$/unidad"
$ws.Range("O772").Value2 = "Provincia de Quillota"
$ws.Range("P772").Value2 = 600
$ws.Range("Q772").Value2 = 1
$ws.Range("R772").Value2 = "Hortaliza"
